$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (164 and 165) following the existing pattern
$ws.Cells.Item(164, 1).Value = 163
$ws.Cells.Item(164, 2).Value = 1
$ws.Cells.Item(164, 3).Value = "2024-06-18 10:13:33"
$ws.Cells.Item(164, 4).Value = 200
$ws.Cells.Item(164, 5).Value = 10

$ws.Cells.Item(165, 1).Value = 164
$ws.Cells.Item(165, 2).Value = 2
$ws.Cells.Item(165, 3).Value = "2024-06-18 10:13:33"
$ws.Cells.Item(165, 4).Value = 200
$ws.Cells.Item(165, 5).Value = 0
